# Updates Leve market-price/profit data across all item-crafting sheets
# (ALC, ARM, CRP, CUL, GSM, LTW, WVR), refreshing currentAveragePrice /
# LevePrice / LeveProfit columns (H-N) per the scheduled market-data sync.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 370.8421
$ws.Cells.Item(15, 9).Value = 370.8421
$ws.Cells.Item(15, 11).Value = 1112.5263
$ws.Cells.Item(15, 13).Value = -943.5263
$ws.Cells.Item(33, 8).Value = 181.5
$ws.Cells.Item(33, 9).Value = 156.2
$ws.Cells.Item(33, 11).Value = 156.2
$ws.Cells.Item(33, 13).Value = 72.80000000000001
$ws.Cells.Item(43, 8).Value = 4606.7407
$ws.Cells.Item(43, 10).Value = 5827.4287
$ws.Cells.Item(43, 12).Value = 5827.4287
$ws.Cells.Item(43, 14).Value = -5965.4287
$ws.Cells.Item(51, 8).Value = 6666.6665
$ws.Cells.Item(51, 9).Value = 5000
$ws.Cells.Item(51, 10).Value = 10000
$ws.Cells.Item(51, 11).Value = 5000
$ws.Cells.Item(51, 12).Value = 10000
$ws.Cells.Item(51, 13).Value = -4516
$ws.Cells.Item(51, 14).Value = -10968
$ws.Cells.Item(63, 8).Value = 40246
$ws.Cells.Item(63, 9).Value = 40246
$ws.Cells.Item(63, 11).Value = 40246
$ws.Cells.Item(63, 13).Value = -39622
$ws.Cells.Item(64, 8).Value = 6420.727
$ws.Cells.Item(64, 9).Value = 3328.625
$ws.Cells.Item(64, 10).Value = 14666.333
$ws.Cells.Item(64, 11).Value = 3328.625
$ws.Cells.Item(64, 12).Value = 14666.333
$ws.Cells.Item(64, 13).Value = -3080.625
$ws.Cells.Item(64, 14).Value = -15162.333
$ws.Cells.Item(66, 8).Value = 40246
$ws.Cells.Item(66, 9).Value = 40246
$ws.Cells.Item(66, 11).Value = 120738
$ws.Cells.Item(66, 13).Value = -117618
$ws.Cells.Item(67, 8).Value = 6420.727
$ws.Cells.Item(67, 9).Value = 3328.625
$ws.Cells.Item(67, 10).Value = 14666.333
$ws.Cells.Item(67, 11).Value = 3328.625
$ws.Cells.Item(67, 12).Value = 14666.333
$ws.Cells.Item(67, 13).Value = -2470.625
$ws.Cells.Item(67, 14).Value = -16382.333
$ws.Cells.Item(70, 8).Value = 2499.5
$ws.Cells.Item(70, 10).Value = 2499.5
$ws.Cells.Item(70, 12).Value = 7498.5
$ws.Cells.Item(70, 14).Value = -8038.5
$ws.Cells.Item(73, 8).Value = 2499.5
$ws.Cells.Item(73, 10).Value = 2499.5
$ws.Cells.Item(73, 12).Value = 7498.5
$ws.Cells.Item(73, 14).Value = -9370.5
$ws.Cells.Item(113, 8).Value = 13663.8
$ws.Cells.Item(113, 9).Value = 18605.834
$ws.Cells.Item(113, 10).Value = 6250.75
$ws.Cells.Item(113, 11).Value = 18605.834
$ws.Cells.Item(113, 12).Value = 6250.75
$ws.Cells.Item(113, 13).Value = -15351.834
$ws.Cells.Item(113, 14).Value = -12758.75
$ws.Cells.Item(141, 8).Value = 651.44446
$ws.Cells.Item(141, 9).Value = 651.44446
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 1954.33338
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = 3225.66662
$ws.Cells.Item(141, 14).ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 3048.2354
$ws.Cells.Item(88, 10).Value = 3336.7856
$ws.Cells.Item(88, 12).Value = 3336.7856
$ws.Cells.Item(88, 14).Value = -4148.7856
$ws.Cells.Item(91, 8).Value = 3048.2354
$ws.Cells.Item(91, 10).Value = 3336.7856
$ws.Cells.Item(91, 12).Value = 3336.7856
$ws.Cells.Item(91, 14).Value = -6144.7856
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 2335.1875
$ws.Cells.Item(7, 9).Value = 1178.4546
$ws.Cells.Item(7, 10).Value = 4880
$ws.Cells.Item(7, 11).Value = 1178.4546
$ws.Cells.Item(7, 12).Value = 4880
$ws.Cells.Item(7, 13).Value = -1065.4546
$ws.Cells.Item(7, 14).Value = -5106
$ws.Cells.Item(22, 8).Value = 656.1667
$ws.Cells.Item(22, 9).Value = 367.58334
$ws.Cells.Item(22, 10).Value = 1233.3334
$ws.Cells.Item(22, 11).Value = 367.58334
$ws.Cells.Item(22, 12).Value = 1233.3334
$ws.Cells.Item(22, 13).Value = -17.58334000000002
$ws.Cells.Item(22, 14).Value = -1933.3334
$ws.Cells.Item(134, 8).Value = 1640.3684
$ws.Cells.Item(134, 9).Value = 1420.7646
$ws.Cells.Item(134, 11).Value = 4262.293799999999
$ws.Cells.Item(134, 13).Value = -1727.293799999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 619
$ws.Cells.Item(5, 9).Value = 619
$ws.Cells.Item(5, 11).Value = 1857
$ws.Cells.Item(5, 13).Value = -1745
$ws.Cells.Item(135, 8).Value = 619
$ws.Cells.Item(135, 9).Value = 619
$ws.Cells.Item(135, 11).Value = 5571
$ws.Cells.Item(135, 13).Value = -3036
$ws.Cells.Item(139, 8).Value = 2977.6667
$ws.Cells.Item(139, 9).Value = 2966.5
$ws.Cells.Item(139, 10).Value = 3000
$ws.Cells.Item(139, 11).Value = 8899.5
$ws.Cells.Item(139, 12).Value = 9000
$ws.Cells.Item(139, 13).Value = -3759.5
$ws.Cells.Item(139, 14).Value = -19280
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 2079
$ws.Cells.Item(46, 9).Value = 2079
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 2079
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -1923
$ws.Cells.Item(46, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 966
$ws.Cells.Item(102, 10).Value = 1331.2
$ws.Cells.Item(102, 12).Value = 1331.2
$ws.Cells.Item(102, 14).Value = -4575.2
$ws.Cells.Item(122, 8).Value = 4110
$ws.Cells.Item(122, 9).Value = 4830.5
$ws.Cells.Item(122, 10).Value = 2669
$ws.Cells.Item(122, 11).Value = 14491.5
$ws.Cells.Item(122, 12).Value = 8007
$ws.Cells.Item(122, 13).Value = -12041.5
$ws.Cells.Item(122, 14).Value = -12907
$ws.Cells.Item(132, 8).Value = 2963.75
$ws.Cells.Item(132, 9).Value = 2963.75
$ws.Cells.Item(132, 11).Value = 8891.25
$ws.Cells.Item(132, 13).Value = -6361.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1500
$ws.Cells.Item(16, 9).Value = 1500
$ws.Cells.Item(16, 11).Value = 1500
$ws.Cells.Item(16, 13).Value = -1330
$ws.Cells.Item(82, 8).Value = 1039.8
$ws.Cells.Item(82, 9).Value = 1099.75
$ws.Cells.Item(82, 10).Value = 800
$ws.Cells.Item(82, 11).Value = 1099.75
$ws.Cells.Item(82, 12).Value = 800
$ws.Cells.Item(82, 13).Value = -738.75
$ws.Cells.Item(82, 14).Value = -1522
$ws.Cells.Item(85, 8).Value = 1039.8
$ws.Cells.Item(85, 9).Value = 1099.75
$ws.Cells.Item(85, 10).Value = 800
$ws.Cells.Item(85, 11).Value = 1099.75
$ws.Cells.Item(85, 12).Value = 800
$ws.Cells.Item(85, 13).Value = 148.25
$ws.Cells.Item(85, 14).Value = -3296
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 50000000
$ws.Cells.Item(49, 10).Value = 50000000
$ws.Cells.Item(49, 12).Value = 50000000
$ws.Cells.Item(49, 14).Value = -50000460
$ws.Cells.Item(62, 8).Value = 2239.4
$ws.Cells.Item(62, 9).Value = 2239.4
$ws.Cells.Item(62, 11).Value = 2239.4
$ws.Cells.Item(62, 13).Value = -1615.4
$ws.Cells.Item(65, 8).Value = 2239.4
$ws.Cells.Item(65, 9).Value = 2239.4
$ws.Cells.Item(65, 11).Value = 11197
$ws.Cells.Item(65, 13).Value = -8077
$ws.Cells.Item(132, 8).Value = 1624.0588
$ws.Cells.Item(132, 9).Value = 1624.0588
$ws.Cells.Item(132, 11).Value = 4872.1764
$ws.Cells.Item(132, 13).Value = -2342.1764
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
Write-Host "Updated Leve profit/price figures across ALC, ARM, CRP, CUL, GSM, LTW, WVR."
